{"js": "// Replace each two-digit-by-two-digit multiplication expression in the\n// document with its new value. Each \"before\" expression is unique in the\n// document, so a literal (non-wildcard, case-sensitive) search-and-replace\n// per pair is sufficient and keeps run formatting (font, size) intact.\nconst replacements = [\n  [\"51\u00d727=1377\", \"73\u00d787=6351\"],\n  [\"14\u00d711=154\", \"37\u00d729=1073\"],\n  [\"99\u00d783=8217\", \"34\u00d794=3196\"],\n  [\"69\u00d711=759\", \"47\u00d738=1786\"],\n  [\"75\u00d752=3900\", \"30\u00d792=2760\"],\n  [\"99\u00d765=6435\", \"17\u00d766=1122\"],\n  [\"51\u00d724=1224\", \"85\u00d730=2550\"],\n  [\"89\u00d756=4984\", \"80\u00d772=5760\"],\n  [\"88\u00d740=3520\", \"69\u00d733=2277\"],\n  [\"34\u00d739=1326\", \"76\u00d717=1292\"],\n  [\"34\u00d798=3332\", \"67\u00d723=1541\"],\n  [\"19\u00d752=988\", \"14\u00d745=630\"],\n  [\"55\u00d793=5115\", \"73\u00d727=1971\"],\n  [\"77\u00d783=6391\", \"51\u00d795=4845\"],\n  [\"34\u00d711=374\", \"51\u00d761=3111\"],\n  [\"32\u00d788=2816\", \"79\u00d753=4187\"],\n  [\"98\u00d797=9506\", \"25\u00d798=2450\"],\n  [\"64\u00d759=3776\", \"43\u00d789=3827\"],\n  [\"52\u00d734=1768\", \"68\u00d792=6256\"],\n  [\"96\u00d746=4416\", \"81\u00d791=7371\"],\n  [\"96\u00d793=8928\", \"31\u00d714=434\"],\n  [\"48\u00d730=1440\", \"46\u00d725=1150\"],\n  [\"39\u00d753=2067\", \"40\u00d760=2400\"],\n  [\"29\u00d775=2175\", \"61\u00d716=976\"],\n  [\"63\u00d781=5103\", \"66\u00d717=1122\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace each two-digit-by-two-digit multiplication expression in the\n# document with its new value. Each \"before\" expression is unique in the\n# document, so a literal (non-wildcard, case-sensitive) find & replace per\n# pair is sufficient and keeps run formatting (font, size) intact.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"51\u00d727=1377\", \"73\u00d787=6351\"),\n  @(\"14\u00d711=154\", \"37\u00d729=1073\"),\n  @(\"99\u00d783=8217\", \"34\u00d794=3196\"),\n  @(\"69\u00d711=759\", \"47\u00d738=1786\"),\n  @(\"75\u00d752=3900\", \"30\u00d792=2760\"),\n  @(\"99\u00d765=6435\", \"17\u00d766=1122\"),\n  @(\"51\u00d724=1224\", \"85\u00d730=2550\"),\n  @(\"89\u00d756=4984\", \"80\u00d772=5760\"),\n  @(\"88\u00d740=3520\", \"69\u00d733=2277\"),\n  @(\"34\u00d739=1326\", \"76\u00d717=1292\"),\n  @(\"34\u00d798=3332\", \"67\u00d723=1541\"),\n  @(\"19\u00d752=988\", \"14\u00d745=630\"),\n  @(\"55\u00d793=5115\", \"73\u00d727=1971\"),\n  @(\"77\u00d783=6391\", \"51\u00d795=4845\"),\n  @(\"34\u00d711=374\", \"51\u00d761=3111\"),\n  @(\"32\u00d788=2816\", \"79\u00d753=4187\"),\n  @(\"98\u00d797=9506\", \"25\u00d798=2450\"),\n  @(\"64\u00d759=3776\", \"43\u00d789=3827\"),\n  @(\"52\u00d734=1768\", \"68\u00d792=6256\"),\n  @(\"96\u00d746=4416\", \"81\u00d791=7371\"),\n  @(\"96\u00d793=8928\", \"31\u00d714=434\"),\n  @(\"48\u00d730=1440\", \"46\u00d725=1150\"),\n  @(\"39\u00d753=2067\", \"40\u00d760=2400\"),\n  @(\"29\u00d775=2175\", \"61\u00d716=976\"),\n  @(\"63\u00d781=5103\", \"66\u00d717=1122\")\n)\n\nforeach ($pair in $pairs) {\n  $old = $pair[0]\n  $new = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $old\n  $find.Replacement.Text = $new\n  $find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null\n}\n"}
